$wb = $excel.ActiveWorkbook

$values = @(0.5228733495041809, -0.40000000000000036, 0.7524446571905798, -0.09109999999999996, 1.5829618029997903, 16.12947350163202, 0.948883713747442)

foreach ($sheetName in @("Test 1", "Test 2")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C11").Value = $values[0]
    $ws.Range("D11").Value = $values[1]
    $ws.Range("E11").Value = $values[2]
    $ws.Range("F11").Value = $values[3]
    $ws.Range("G11").Value = $values[4]
    $ws.Range("H11").Value = $values[5]
    $ws.Range("I11").Value = $values[6]
}
